# Auto-generated edit script: updates Behemoth_Profits price/profit data across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ALC_values = @{
    "H11" = 1403.5
    "I11" = 1403.5
    "K11" = 1403.5
    "M11" = -1263.5
    "H31" = 453.77777
    "J31" = 445.66666
    "L31" = 1336.99998
    "N31" = -1796.99998
    "H38" = 406.33334
    "I38" = 285.8
    "K38" = 857.4000000000001
    "M38" = -485.4000000000001
    "H39" = 278.89285
    "I39" = 240
    "J39" = 300.5
    "K39" = 720
    "L39" = 901.5
    "M39" = -424
    "N39" = -1493.5
    "H41" = 316.85715
    "I41" = 350.25
    "K41" = 350.25
    "M41" = 89.75
    "H42" = 306
    "I42" = 248
    "K42" = 744
    "M42" = -514
    "H51" = 96269.84
    "I51" = 104050.9
    "K51" = 104050.9
    "M51" = -103566.9
    "H64" = 6227.5
    "I64" = 4900
    "K64" = 4900
    "M64" = -4652
    "H67" = 6227.5
    "I67" = 4900
    "K67" = 4900
    "M67" = -4042
    "H75" = 0
    "J75" = 0
    "L75" = 0
    "H76" = 6218.5557
    "I76" = 3991.75
    "J76" = 8000
    "K76" = 3991.75
    "L76" = 8000
    "M76" = -3676.75
    "N76" = -8630
    "H78" = 0
    "J78" = 0
    "L78" = 0
    "H79" = 6218.5557
    "I79" = 3991.75
    "J79" = 8000
    "K79" = 3991.75
    "L79" = 8000
    "M79" = -2899.75
    "N79" = -10184
    "H86" = 2226821
    "J86" = 5371.4287
    "L86" = 5371.4287
    "N86" = -7617.4287
    "H87" = 88419.336
    "J87" = 88419.336
    "L87" = 88419.336
    "N87" = -90915.336
    "H89" = 2226821
    "J89" = 5371.4287
    "L89" = 26857.1435
    "N89" = -38089.14350000001
    "H90" = 88419.336
    "J90" = 88419.336
    "L90" = 265258.008
    "N90" = -277738.008
    "H97" = 1100
    "J97" = 1100
    "L97" = 3300
    "N97" = -4292
    "H112" = 1883.3334
    "I112" = 250
    "K112" = 750
    "M112" = 358
    "H113" = 166667000
    "I113" = 50000500
    "J113" = 400000000
    "K113" = 50000500
    "L113" = 400000000
    "M113" = -49997246
    "N113" = -400006508
    "H116" = 5144.7915
    "I116" = 4332.278
    "K116" = 4332.278
    "M116" = -890.2780000000002
    "H125" = 2761.5833
    "I125" = 1830
    "K125" = 16470
    "M125" = -14010
}
foreach ($key in $ALC_values.Keys) {
    $ws.Range($key).Value = $ALC_values[$key]
}
$ALC_clears = @("N75", "N78")
foreach ($key in $ALC_clears) {
    $ws.Range($key).ClearContents()
}

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ARM_values = @{
    "H32" = 13899993
    "I32" = 16675286
    "K32" = 16675286
    "M32" = -16674999
    "H45" = 55558644
    "I45" = 100002150
    "J45" = 4258
    "K45" = 100002150
    "L45" = 4258
    "M45" = -100001773
    "N45" = -5012
    "H61" = 32614074
    "I61" = 25004824
    "K61" = 25004824
    "M61" = -25004612
    "H69" = 0
    "J69" = 0
    "L69" = 0
    "H72" = 0
    "J72" = 0
    "L72" = 0
    "H74" = 12385076
    "I74" = 16667733
    "K74" = 16667733
    "M74" = -16666859
    "H77" = 12385076
    "I77" = 16667733
    "K77" = 83338665
    "M77" = -83334297
    "H106" = 52240
    "J106" = 52240
    "L106" = 52240
    "N106" = -54764
    "H122" = 4999.5
    "I122" = 4999
    "K122" = 14997
    "M122" = -12547
    "H136" = 32614074
    "I136" = 25004824
    "K136" = 75014472
    "M136" = -75011922
}
foreach ($key in $ARM_values.Keys) {
    $ws.Range($key).Value = $ARM_values[$key]
}
$ARM_clears = @("N69", "N72")
foreach ($key in $ARM_clears) {
    $ws.Range($key).ClearContents()
}

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$BSM_values = @{
    "H20" = 4332.9287
    "I20" = 3972.1667
    "K20" = 3972.1667
    "M20" = -3725.1667
    "H126" = 125000.5
    "J126" = 125000.5
    "L126" = 125000.5
    "N126" = -134880.5
}
foreach ($key in $BSM_values.Keys) {
    $ws.Range($key).Value = $BSM_values[$key]
}

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$CRP_values = @{
    "H2" = 110000
    "J2" = 0
    "L2" = 0
    "H7" = 296.14285
    "I7" = 215
    "K7" = 215
    "M7" = -102
    "H16" = 1088.7142
    "I16" = 774
    "K16" = 774
    "M16" = -487
    "H80" = 89999.664
    "J80" = 89999.664
    "L80" = 89999.664
    "N80" = -92245.664
    "H83" = 89999.664
    "J83" = 89999.664
    "L83" = 269998.992
    "N83" = -281230.992
    "H113" = 1088.7142
    "I113" = 774
    "K113" = 774
    "M113" = 1396
    "H132" = 4079.9092
    "I132" = 3688.1
    "K132" = 11064.3
    "M132" = -8534.299999999999
    "H134" = 3851.5
    "I134" = 2660.8572
    "J134" = 5518.4
    "K134" = 7982.571599999999
    "L134" = 16555.2
    "M134" = -5447.571599999999
    "N134" = -21625.2
    "H140" = 0
    "J140" = 0
    "L140" = 0
    "H141" = 233878.69
    "J141" = 242870.6
    "L141" = 242870.6
    "N141" = -253230.6
}
foreach ($key in $CRP_values.Keys) {
    $ws.Range($key).Value = $CRP_values[$key]
}
$CRP_clears = @("N2", "N140")
foreach ($key in $CRP_clears) {
    $ws.Range($key).ClearContents()
}

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$CUL_values = @{
    "H4" = 32299952
    "I4" = 20166932
    "K4" = 60500796
    "M4" = -60500684
    "H14" = 622.75
    "I14" = 622.75
    "K14" = 1868.25
    "M14" = -1695.25
    "H69" = 2775
    "J69" = 3106.25
    "L69" = 9318.75
    "N69" = -10940.75
    "H72" = 2775
    "J72" = 3106.25
    "L72" = 27956.25
    "N72" = -36068.25
    "H74" = 13208.333
    "J74" = 13208.333
    "L74" = 39624.999
    "N74" = -41746.999
    "H77" = 13208.333
    "J77" = 13208.333
    "L77" = 118874.997
    "N77" = -129482.997
    "H107" = 599.931
    "J107" = 831.3
    "L107" = 2493.9
    "N107" = -6333.9
    "H109" = 1588.6666
    "I109" = 1588.6666
    "K109" = 4765.9998
    "M109" = -3725.9998
    "H117" = 963.75
    "I117" = 287.5
    "K117" = 862.5
    "M117" = 2579.5
    "H118" = 6199.7144
    "J118" = 9125
    "L118" = 27375
    "N118" = -29861
    "H119" = 7999.5557
    "I119" = 2399.2
    "K119" = 7197.599999999999
    "M119" = -2359.599999999999
    "H132" = 1727.2727
    "J132" = 2143.1428
    "L132" = 19288.2852
    "N132" = -24348.2852
    "H137" = 5341
    "J137" = 7499.5
    "L137" = 22498.5
    "N137" = -32698.5
}
foreach ($key in $CUL_values.Keys) {
    $ws.Range($key).Value = $CUL_values[$key]
}

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$GSM_values = @{
    "H15" = 45000
    "J15" = 45000
    "L15" = 45000
    "N15" = -45576
    "H70" = 10992.412
    "I70" = 12287.3
    "J70" = 9142.571
    "K70" = 12287.3
    "L70" = 9142.571
    "M70" = -12017.3
    "N70" = -9682.571
    "H73" = 10992.412
    "I73" = 12287.3
    "J73" = 9142.571
    "K73" = 12287.3
    "L73" = 9142.571
    "M73" = -11351.3
    "N73" = -11014.571
    "H80" = 3840.2666
    "I80" = 3291.5
    "K80" = 3291.5
    "M80" = -2293.5
    "H81" = 45000
    "J81" = 45000
    "L81" = 45000
    "N81" = -46996
    "H83" = 3840.2666
    "I83" = 3291.5
    "K83" = 16457.5
    "M83" = -11465.5
    "H84" = 45000
    "J84" = 45000
    "L84" = 135000
    "N84" = -144984
    "H98" = 15966.286
    "J98" = 18485.666
    "L98" = 18485.666
    "N98" = -24475.666
    "H102" = 1872.6111
    "I102" = 1323.7407
    "J102" = 3519.2222
    "K102" = 1323.7407
    "L102" = 3519.2222
    "M102" = 298.2592999999999
    "N102" = -6763.2222
    "H104" = 99980
    "J104" = 99980
    "L104" = 99980
    "N104" = -106968
    "H113" = 3087.682
    "I113" = 2034.909
    "K113" = 2034.909
    "M113" = 135.0909999999999
    "H122" = 2469.1667
    "I122" = 1851.75
    "K122" = 5555.25
    "M122" = -3105.25
    "H132" = 90920776
    "J132" = 20443
    "L132" = 61329
    "N132" = -66389
}
foreach ($key in $GSM_values.Keys) {
    $ws.Range($key).Value = $GSM_values[$key]
}

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$LTW_values = @{
    "H16" = 1878.6364
    "I16" = 1878.6364
    "K16" = 1878.6364
    "M16" = -1708.6364
    "H22" = 1486
    "I22" = 1380
    "J22" = 1751
    "K22" = 1380
    "L22" = 1751
    "M22" = -1085
    "N22" = -2341
    "H27" = 1486
    "I27" = 1380
    "J27" = 1751
    "K27" = 1380
    "L27" = 1751
    "M27" = -1273
    "N27" = -1965
    "H50" = 33999.5
    "I50" = 28000
    "J50" = 39999
    "K50" = 28000
    "L50" = 39999
    "M50" = -27363
    "N50" = -41273
    "H55" = 30303624
    "I55" = 38462100
    "J55" = 709.8570999999999
    "K55" = 38462100
    "L55" = 709.8570999999999
    "M55" = -38461927
    "N55" = -1055.8571
    "H68" = 4000
    "J68" = 0
    "L68" = 0
    "H71" = 4000
    "J71" = 0
    "L71" = 0
    "H96" = 95750
    "J96" = 95750
    "L96" = 95750
    "N96" = -101242
    "H109" = 98291.664
    "J109" = 98291.664
    "L109" = 98291.664
    "N109" = -101065.664
    "H122" = 5859.4165
    "I122" = 5480.4614
    "J122" = 6307.273
    "K122" = 16441.3842
    "L122" = 18921.819
    "M122" = -13991.3842
    "N122" = -23821.819
    "H136" = 78137.89
    "I136" = 10769.728
    "K136" = 32309.184
    "M136" = -29759.184
}
foreach ($key in $LTW_values.Keys) {
    $ws.Range($key).Value = $LTW_values[$key]
}
$LTW_clears = @("N68", "N71")
foreach ($key in $LTW_clears) {
    $ws.Range($key).ClearContents()
}

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$WVR_values = @{
    "H4" = 1344375.8
    "I4" = 2127000
    "J4" = 40002
    "K4" = 2127000
    "L4" = 40002
    "M4" = -2126887
    "N4" = -40228
    "H5" = 4302000
    "I5" = 377499.75
    "K5" = 377499.75
    "M5" = -377387.75
    "H10" = 0
    "I10" = 0
    "K10" = 0
    "H11" = 0
    "I11" = 0
    "J11" = 0
    "K11" = 0
    "L11" = 0
    "H17" = 0
    "I17" = 0
    "K17" = 0
    "H22" = 1115
    "J22" = 1115
    "L22" = 1115
    "N22" = -1701
    "H23" = 10635722
    "I23" = 652.25
    "K23" = 652.25
    "M23" = -423.25
    "H26" = 0
    "I26" = 0
    "K26" = 0
    "H31" = 61307.6
    "J31" = 73012.664
    "L31" = 73012.664
    "N31" = -73708.664
    "H75" = 9443155
    "J75" = 10785034
    "L75" = 10785034
    "N75" = -10786906
    "H78" = 9443155
    "J78" = 10785034
    "L78" = 32355102
    "N78" = -32364462
    "H81" = 2222.3333
    "I81" = 2083.5
    "J81" = 2500
    "K81" = 4167
    "L81" = 5000
    "M81" = -3106
    "N81" = -7122
    "H84" = 2222.3333
    "I84" = 2083.5
    "J84" = 2500
    "K84" = 20835
    "L84" = 25000
    "M84" = -15531
    "N84" = -35608
    "H107" = 833
    "I107" = 949.8
    "K107" = 2849.4
    "M107" = -929.3999999999996
    "H109" = 104930
    "J109" = 104930
    "L109" = 104930
    "N109" = -107704
    "H122" = 3950
    "I122" = 3285
    "J122" = 5501.6665
    "K122" = 9855
    "L122" = 16504.9995
    "M122" = -7405
    "N122" = -21404.9995
    "H132" = 672400.2
    "I132" = 5199.9
    "K132" = 15599.7
    "M132" = -13069.7
    "H135" = 0
    "J135" = 0
    "L135" = 0
}
foreach ($key in $WVR_values.Keys) {
    $ws.Range($key).Value = $WVR_values[$key]
}
$WVR_clears = @("M10", "M11", "N11", "M17", "M26", "N135")
foreach ($key in $WVR_clears) {
    $ws.Range($key).ClearContents()
}
